$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 (Saurabh Tiwary / Mumbai Indians) had their runs/balls/fours
# values swapped between each other ("updated activity till excel form").
# Keep the cells stored as text (matching the rest of the sheet) by
# formatting the range as text before writing the new values.
$rng = $ws.Range("C5:E6")
$rng.NumberFormat = "@"

$ws.Range("C5").Value = "42"
$ws.Range("D5").Value = "31"
$ws.Range("E5").Value = "3"

$ws.Range("C6").Value = "21"
$ws.Range("D6").Value = "13"
$ws.Range("E6").Value = "1"
